# "custom colours.xlsx" re-style:
#   - insert a new row 1 for a "warm-yellow-pale" swatch (pushes the
#     existing 10 rows down to rows 2-11)
#   - fill its A1 swatch cell with a new custom colour (FFF7E5BC), which
#     mints a new fill + cellXfs entry, same as picking a custom colour in
#     the Fill Color picker
#   - label it in column B via a new shared string
#   - leave the sheet selection on C13 and the page orientation set to
#     portrait, matching the saved session state

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing 10 colour rows down and create the new first row.
$ws.Rows("1:1").Insert()

# New swatch label for row 1.
$ws.Range("B1").Value = "warm-yellow-pale"

# Custom fill colour FFF7E5BC (R=247,G=229,B=188) for the new swatch cell.
# Excel COM colours are packed as 0xBBGGRR.
$paleYellow = 188 * 65536 + 229 * 256 + 247
$ws.Range("A1").Interior.Color = $paleYellow

# Restore the saved selection/view state.
$null = $ws.Range("C13").Select()
$ws.PageSetup.Orientation = 1
